# Apply the refreshed coin market snapshot to the active worksheet.
#
# Columns D (Price), E (Volume 1h %) and G (Hora) hold numeric-looking
# text in the source workbook (e.g. "275.41", "-1.32%", "10"). They are
# assigned here with a leading apostrophe ( "'" + value ) so Excel keeps
# storing them as literal text instead of converting them to numbers,
# which preserves formatting such as trailing zeros (e.g. "6.910").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + '275.41'
$ws.Range("E2").Value = "'" + '-1.32%'
$ws.Range("G2").Value = "'" + '10'

# Row 3
$ws.Range("D3").Value = "'" + '26.64'
$ws.Range("E3").Value = "'" + '-2.49%'
$ws.Range("G3").Value = "'" + '10'

# Row 4
$ws.Range("D4").Value = "'" + '4.878'
$ws.Range("E4").Value = "'" + '1.57%'
$ws.Range("G4").Value = "'" + '10'

# Row 5
$ws.Range("D5").Value = "'" + '0.06343'
$ws.Range("E5").Value = "'" + '0.36%'
$ws.Range("G5").Value = "'" + '10'

# Row 6
$ws.Range("D6").Value = "'" + '6.910'
$ws.Range("E6").Value = "'" + '-0.23%'
$ws.Range("G6").Value = "'" + '10'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'" + '1.262'
$ws.Range("E7").Value = "'" + '33.27%'
$ws.Range("G7").Value = "'" + '10'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'" + '0.8736'
$ws.Range("E8").Value = "'" + '-0.67%'
$ws.Range("G8").Value = "'" + '10'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = "'" + '0.1533'
$ws.Range("E9").Value = "'" + '4.16%'
$ws.Range("G9").Value = "'" + '10'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = "'" + '0.05029'
$ws.Range("E10").Value = "'" + '-1.04%'
$ws.Range("G10").Value = "'" + '10'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'" + '0.07419'
$ws.Range("E11").Value = "'" + '1.15%'
$ws.Range("G11").Value = "'" + '10'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'" + '0.02973'
$ws.Range("E12").Value = "'" + '-5.69%'
$ws.Range("G12").Value = "'" + '10'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'" + '0.09056'
$ws.Range("E13").Value = "'" + '-0.11%'
$ws.Range("G13").Value = "'" + '10'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'" + '0.001581'
$ws.Range("E14").Value = "'" + '2.09%'
$ws.Range("G14").Value = "'" + '10'

# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'" + '0.0006301'
$ws.Range("E15").Value = "'" + '0.49%'
$ws.Range("G15").Value = "'" + '10'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'" + '0.006026'
$ws.Range("E16").Value = "'" + '2.24%'
$ws.Range("G16").Value = "'" + '10'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'" + '3.447'
$ws.Range("E17").Value = "'" + '0.14%'
$ws.Range("G17").Value = "'" + '10'

# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = "'" + '3.317'
$ws.Range("E18").Value = "'" + '0.25%'
$ws.Range("G18").Value = "'" + '10'

# Row 19
$ws.Range("D19").Value = "'" + '2.271'
$ws.Range("E19").Value = "'" + '-0.88%'
$ws.Range("G19").Value = "'" + '10'

# Row 20
$ws.Range("E20").Value = "'" + '-1.08%'
$ws.Range("G20").Value = "'" + '10'

# Row 21
$ws.Range("D21").Value = "'" + '0.1324'
$ws.Range("E21").Value = "'" + '1.00%'
$ws.Range("G21").Value = "'" + '10'

# Row 22
$ws.Range("D22").Value = "'" + '3.910'
$ws.Range("E22").Value = "'" + '0.74%'
$ws.Range("G22").Value = "'" + '10'

# Row 23
$ws.Range("D23").Value = "'" + '0.04377'
$ws.Range("E23").Value = "'" + '0.87%'
$ws.Range("G23").Value = "'" + '10'

# Row 24
$ws.Range("D24").Value = "'" + '0.001172'
$ws.Range("E24").Value = "'" + '-1.04%'
$ws.Range("G24").Value = "'" + '10'

# Row 25
$ws.Range("D25").Value = "'" + '0.004210'
$ws.Range("E25").Value = "'" + '-1.66%'
$ws.Range("G25").Value = "'" + '10'

# Row 26
$ws.Range("D26").Value = "'" + '0.0001201'
$ws.Range("E26").Value = "'" + '0.11%'
$ws.Range("G26").Value = "'" + '10'

# Row 27
$ws.Range("D27").Value = "'" + '0.0001680'
$ws.Range("E27").Value = "'" + '-0.58%'
$ws.Range("G27").Value = "'" + '10'

# Row 28
$ws.Range("G28").Value = "'" + '10'

# Row 29
$ws.Range("G29").Value = "'" + '10'

# Row 30
$ws.Range("G30").Value = "'" + '10'

# Row 31
$ws.Range("G31").Value = "'" + '10'

# Row 32
$ws.Range("G32").Value = "'" + '10'

# Row 33
$ws.Range("G33").Value = "'" + '10'

# Row 34
$ws.Range("G34").Value = "'" + '10'

# Row 35
$ws.Range("G35").Value = "'" + '10'

# Row 36
$ws.Range("G36").Value = "'" + '10'

# Row 37
$ws.Range("G37").Value = "'" + '10'

# Row 38
$ws.Range("G38").Value = "'" + '10'

# Row 39
$ws.Range("G39").Value = "'" + '10'

# Row 40
$ws.Range("D40").Value = "'" + '0.04107'
$ws.Range("E40").Value = "'" + '0.83%'
$ws.Range("G40").Value = "'" + '10'

# Row 41
$ws.Range("D41").Value = "'" + '0.006980'
$ws.Range("E41").Value = "'" + '6.14%'
$ws.Range("G41").Value = "'" + '10'

# Row 42
$ws.Range("E42").Value = "'" + '1.07%'
$ws.Range("G42").Value = "'" + '10'

# Row 43
$ws.Range("E43").Value = "'" + '-2.53%'
$ws.Range("G43").Value = "'" + '10'

# Row 44
$ws.Range("D44").Value = "'" + '0.01080'
$ws.Range("E44").Value = "'" + '-16.80%'
$ws.Range("G44").Value = "'" + '10'

# Row 45
$ws.Range("D45").Value = "'" + '0.00005299'
$ws.Range("E45").Value = "'" + '1.67%'
$ws.Range("G45").Value = "'" + '10'

# Row 46
$ws.Range("D46").Value = "'" + '0.02102'
$ws.Range("E46").Value = "'" + '-6.62%'
$ws.Range("G46").Value = "'" + '10'

# Row 47
$ws.Range("G47").Value = "'" + '10'

# Row 48
$ws.Range("G48").Value = "'" + '10'

# Row 49
$ws.Range("G49").Value = "'" + '10'

# Row 50
$ws.Range("G50").Value = "'" + '10'

# Row 51
$ws.Range("G51").Value = "'" + '10'
